# Update the "Yoff" column (K) for each muscle row so that it reflects
# -1/10 of the "Am" column (I) value, instead of the placeholder -1E-3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 39; $row++) {
    $amValue = $ws.Cells.Item($row, 9).Value()   # Column I = Am
    $ws.Cells.Item($row, 11).Value = (-1 * $amValue) / 10   # Column K = Yoff
}
